# Update cryptocurrency price/volume data per the Dec 5 2023 GitHub Actions refresh.
# Cells whose new text is a bare number (e.g. "231.57") are written with a leading
# apostrophe (like typing '231.57 in Excel) so they stay text, matching the source
# column's formatting (prices are text strings like "41.774.58", not numeric values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.774.58"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.230.20"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'231.57"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'0.619"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("D7").Value = "'60.47"
$ws.Range("E7").Value = "  -5.81%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'58.24"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "'0.0904"
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "2.562.01"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'15.77"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'22.87"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "2.252.44"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").Value = "41.736.72"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'72.34"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "'248.75"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'9.76"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").Value = "'169.59"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "'0.143"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "'2.65"
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").Value = "'5.06"
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "'6.59"
$ws.Range("E37").Value = "  -7.41%  "
$ws.Range("D38").Value = "'3.63"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("B40").Value = "TerraClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D40").Value = "'0.000238"
$ws.Range("E40").Value = "  +14.67%  "
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +3.30%  "
$ws.Range("D43").Value = "'8.53"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "'98.68"
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").Value = "'4.40"
$ws.Range("E47").Value = "  -8.33%  "
$ws.Range("D48").Value = "1.470.32"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  +11.12%  "
